$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Bolzano-Bozen" and "Trento" rows (rows 7 and 8), shifting
# everything below up by two rows. This also drops the two now-unused
# shared strings ("Bolzano-Bozen ", "Trento ") from the shared string table.
$ws.Range("A7:G8").EntireRow.Delete()

# Reflect the reviewer's final on-screen state: zoomed in, selection left
# on the row that used to hold "Emilia-Romagna" (now a few rows higher).
$ws.Select() | Out-Null
$excel.ActiveWindow.Zoom = 145
$ws.Range("J15").Select() | Out-Null
